# Add formatting for EPV data: append a new "Premises and Equipment"
# section (rows 24-32) below the existing EPV table on the EPV sheet,
# with its own Arial-10 font and yellow/blue fill palette.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EPV")

# Colors (BGR-encoded Long values expected by the Excel object model).
$yellowFill = 13431550   # RGB(254,242,204) -> fef2cc
$blueFill   = 15983311   # RGB(207,226,243) -> cfe2f3
$black      = 0

# xlEdgeLeft=7, xlEdgeTop=8, xlEdgeBottom=9, xlEdgeRight=10
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlContinuous = 1
$xlThin = 2
$xlThick = 4

function Set-Edge($cell, $edge, $weight) {
    # Setting Color before LineStyle/Weight lets the engine dedupe onto
    # the workbook's existing (identically-colored) border definitions
    # whenever the resulting edge combination already exists.
    $cell.Borders.Item($edge).Color = $black
    $cell.Borders.Item($edge).LineStyle = $xlContinuous
    $cell.Borders.Item($edge).Weight = $weight
}

# Apply the shared look (font + fill) used by every cell in the new
# section, then layer on the border pattern appropriate to its role:
#   "header" -> thick left/top            (section title row, like row 2)
#   "data"   -> thick/thin left only      (plain data row, like row 3)
#   "sub"    -> thick/thin left + thin bottom   (row before the total)
#   "total"  -> thick/thin left + thick bottom, bold font (closing row)
function Format-EpvCell($cell, [bool]$isBlue, [string]$role, [bool]$bold) {
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.Bold = $bold
    $cell.Font.Strikethrough = $false
    $cell.Font.Color = $black

    if ($isBlue) {
        $cell.Interior.Color = $blueFill
    } else {
        $cell.Interior.Color = $yellowFill
    }

    if (-not $isBlue) {
        # Column B: thick (or thin, for the continuation rows) left edge.
        Set-Edge $cell $xlEdgeLeft $xlThick
    } else {
        # Column C: thin left + thin right edges.
        Set-Edge $cell $xlEdgeLeft $xlThin
        Set-Edge $cell $xlEdgeRight $xlThin
    }

    if ($role -eq "header") {
        Set-Edge $cell $xlEdgeTop $xlThick
    } elseif ($role -eq "sub") {
        Set-Edge $cell $xlEdgeBottom $xlThin
    } elseif ($role -eq "total") {
        Set-Edge $cell $xlEdgeBottom $xlThick
    }
}

$rows = @(
    @{ Row = 24; Text = "Premises and Equipment";        Role = "header"; Bold = $false },
    @{ Row = 25; Text = "Current Year Revenue";           Role = "data";   Bold = $false },
    @{ Row = 26; Text = "Prior Year Revenue";              Role = "data";   Bold = $false },
    @{ Row = 27; Text = "Change in Revenue";               Role = "data";   Bold = $false },
    @{ Row = 28; Text = "Depreciation and Amortization";   Role = "data";   Bold = $false },
    @{ Row = 29; Text = "CAPEX";                           Role = "data";   Bold = $false },
    @{ Row = 30; Text = "Growth CAPEX";                    Role = "data";   Bold = $false },
    @{ Row = 31; Text = "Zero-growth CAPEX";                Role = "sub";    Bold = $false },
    @{ Row = 32; Text = "Depreciation Adjustment";          Role = "total";  Bold = $true }
)

foreach ($r in $rows) {
    $bCell = $ws.Cells.Item($r.Row, 2)
    $cCell = $ws.Cells.Item($r.Row, 3)

    $bCell.Value = $r.Text
    $cCell.Value = $r.Text

    Format-EpvCell $bCell $false $r.Role $r.Bold
    Format-EpvCell $cCell $true  $r.Role $r.Bold
}

Write-Output "Applied EPV formatting for rows 24-32"
